# Apply the "added bbc article classifier" edit to Sheet1 of the workbook.
# This fills in the worked Naive Bayes example (word counts, probabilities,
# Laplace-smoothed probabilities, and the final product rows) and adds a
# couple of new labels, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New header / label cells -------------------------------------------------
$ws.Range("F2").Value = "tag = spam"
$ws.Range("F3").Value = "tag = ham"

$ws.Range("E6").Value = "Spam"
$ws.Range("A7").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D7").Value = "P(C|spam)*P(you|spam)…"

# --- Totals used throughout the sheet -----------------------------------------
$ws.Range("B8").Value = 7
$ws.Range("B9").Value = 14
$ws.Range("D9").Value = " "
$ws.Range("B10").Value = 17

# --- Row 13 header: move "Count"/"Probability" one column to the right --------
$ws.Range("D13").ClearContents()
$ws.Range("E13").Value = "Count"
$ws.Range("F13").Value = "Probability"

# --- Step2 table: word counts and probabilities (rows 14-18) ------------------
$ws.Range("B14").Value = 0
$ws.Range("C14").Formula = "=B14/`$B`$8"

$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 1
$ws.Range("B17").Value = 2
# C15 picks up the plain "0.0000" style (index 7), same as its neighbours,
# instead of the comma-style variant (index 8) it had before.
$ws.Range("C15").ClearFormats()
$ws.Range("C16").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C15:C17").Formula = "=B15/`$B`$8"

$ws.Range("B18").Value = 2
$ws.Range("C18").Formula = "=B18/5"

$ws.Range("E14").Value = 0
$ws.Range("F14").Formula = "=E14/`$B`$9"

$ws.Range("E15").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F15:F17").Formula = "=E15/`$B`$9"

$ws.Range("E18").Value = 3
$ws.Range("F18").Formula = "=E18/5"

# --- New row 19: product of the five probabilities -----------------------------
$ws.Range("C19").Formula = "=C14*C15*C16*C17*C18"
$ws.Range("F19").Formula = "=F14*F15*F16*F17*F18"
$ws.Range("C30").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("F19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("A31").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step3 Laplace-smoothed table (rows 26-29) ---------------------------------
$ws.Range("B26").Value = 0
$ws.Range("B27").Value = 1
$ws.Range("B28").Value = 1
$ws.Range("B29").Value = 2

$ws.Range("E26").Value = 0
$ws.Range("E27").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("E29").Value = 0

# --- Row 30: final Laplace counts and probabilities -----------------------------
$ws.Range("B30").ClearFormats()
$ws.Range("B30").Value = 2
$ws.Range("E30").ClearFormats()
$ws.Range("E30").Value = 3
$ws.Range("F30").Formula = "=E30/5"

Write-Output "done"
